$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert 7 new rows after row 16 for the new "cholera_cases" group.
$ws.Rows.Item(17).Resize(7).Insert()

# Row insertion only propagates formatting in column A (the only populated
# column in the template row 16), so stamp the default "style 1" formatting
# (matches every other data cell, e.g. A15/B15) across columns B:D and F of
# the new question rows (17:22) before writing values. Column E is left
# untouched except on row 20, which is the only row that actually uses it
# (the relevant-expression cell). Row 23 (end group) only needs column A,
# which already inherited the right style from the insert.
$defaultFmt = $ws.Cells.Item(15,1)
$defaultFmt.Copy()
$ws.Range("B17:D22").PasteSpecial(-4122)
$defaultFmt.Copy()
$ws.Range("F17:F22").PasteSpecial(-4122)
$defaultFmt.Copy()
$ws.Cells.Item(20,5).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 18's label/required/appearance cells use the alternate "style 2"
# formatting (matches the existing C15/H8 cells). PasteSpecial only
# applies to the first area of a multi-area range here, so paste each
# target cell individually.
$altFmt = $ws.Cells.Item(15,3)
$altFmt.Copy()
$ws.Cells.Item(18,3).PasteSpecial(-4122)
$altFmt.Copy()
$ws.Cells.Item(18,4).PasteSpecial(-4122)
$altFmt.Copy()
$ws.Cells.Item(18,6).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 17: begin group cholera_cases
$ws.Cells.Item(17,1).Value = "begin group"
$ws.Cells.Item(17,2).Value = "cholera_cases"
$ws.Cells.Item(17,3).Value = "Cholera Cases"
# D17 / F17 stay touched-but-empty cells.

# Row 18: num_screened
$ws.Cells.Item(18,1).Value = "integer"
$ws.Cells.Item(18,2).Value = "num_screened"
$ws.Cells.Item(18,3).Value = "Number Of Persons Screened For Cholera"
$ws.Cells.Item(18,4).Value = "yes"
$ws.Cells.Item(18,6).Value = "numbers"

# Row 19: num_referred
$ws.Cells.Item(19,1).Value = "integer"
$ws.Cells.Item(19,2).Value = "num_referred"
$ws.Cells.Item(19,3).Value = "Number Of Presumptive Cholera Persons Referred For Diagnosis"
$ws.Cells.Item(19,4).Value = "yes"
$ws.Cells.Item(19,6).Value = "numbers"

# Row 20: num_referred_reached (adds a relevant formula in column E)
$ws.Cells.Item(20,1).Value = "integer"
$ws.Cells.Item(20,2).Value = "num_referred_reached"
$ws.Cells.Item(20,3).Value = "Number Of Referred Persons Who Reached Health Facility"
$ws.Cells.Item(20,4).Value = "yes"
$ws.Cells.Item(20,5).Value = '${num_referred} >0'
$ws.Cells.Item(20,6).Value = "numbers"

# Row 21: num_confirmed_cases
$ws.Cells.Item(21,1).Value = "integer"
$ws.Cells.Item(21,2).Value = "num_confirmed_cases"
$ws.Cells.Item(21,3).Value = "Number Of Confirmed Cholera Cases At Health Facility"
$ws.Cells.Item(21,4).Value = "yes"
$ws.Cells.Item(21,6).Value = "numbers"

# Row 22: num_deaths
$ws.Cells.Item(22,1).Value = "integer"
$ws.Cells.Item(22,2).Value = "num_deaths"
$ws.Cells.Item(22,3).Value = "Number Of Deaths Due To Cholera In The Month"
$ws.Cells.Item(22,4).Value = "yes"
$ws.Cells.Item(22,6).Value = "numbers"

# Row 23: end group
$ws.Cells.Item(23,1).Value = "end group"

# New column E width for the relevant-expression column we just populated
# (ColumnWidth is character-width units; the saved XML "width" runs
# ColumnWidth + 5/6, so back the value off to land on exactly 22.0).
$ws.Columns.Item(5).ColumnWidth = 21.166666666666668
